# Add a "Skill Description" column (full/friendly names for each SkillCode)
# right after the existing "SkillCode" column, shifting the old
# "SFIA Level" / "Keycode" / "Description" columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map from the short SkillCode / category code (column A) to its full name.
$fullNames = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "REQM"       = "Requirements definition and management"
    "BPTS"       = "Acceptance testing"
    "DTAN"       = "Data modelling and design"
    "MADE"       = "MADE"
    "BSMO"       = "Business modelling"
}

# Insert a new, empty column before column B - this pushes the old
# B (SFIA Level), C (Keycode) and D (Description) columns to C, D, E.
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Skill Description"

# Fill in the full name for every data row, based on the SkillCode in column A.
$lastRow = 25
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value()
    if ($code -ne $null -and $code -ne "") {
        $ws.Cells.Item($r, 2).Value = $fullNames[$code]
    }
}
